$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "71.728.76"
$ws.Range("E2").Value = "  +0.56%  "

$ws.Range("D3").Value = "3.986.24"
$ws.Range("E3").Value = "  -0.43%  "

$ws.Range("E4").Value = "  -0.01%  "

$ws.Range("D5").Value = "527.60"
$ws.Range("E5").Value = "  -0.19%  "

$ws.Range("D6").Value = "149.35"
$ws.Range("E6").Value = "  +0.22%  "

$ws.Range("D7").Value = "0.689"
$ws.Range("E7").Value = "  +10.81%  "

$ws.Range("E8").Value = "  -0.01%  "

$ws.Range("D9").Value = "0.739"
$ws.Range("E9").Value = "  +0.32%  "

$ws.Range("D10").Value = "0.170"
$ws.Range("E10").Value = "  -3.81%  "

$ws.Range("D11").Value = "0.0000322"
$ws.Range("E11").Value = "  -5.74%  "

$ws.Range("D12").Value = "46.73"
$ws.Range("E12").Value = "  +8.16%  "

$ws.Range("D13").Value = "4.635.19"
$ws.Range("E13").Value = "  -0.02%  "

$ws.Range("D14").Value = "10.52"
$ws.Range("E14").Value = "  -1.27%  "

$ws.Range("D15").Value = "3.995.05"
$ws.Range("E15").Value = "  -0.29%  "

$ws.Range("D16").Value = "13.84"
$ws.Range("E16").Value = "  -3.58%  "

$ws.Range("D17").Value = "20.26"
$ws.Range("E17").Value = "  -5.06%  "

$ws.Range("E18").Value = "  -1.18%  "

$ws.Range("D19").Value = "1.17"
$ws.Range("E19").Value = "  -4.43%  "

$ws.Range("D20").Value = "71.563.08"
$ws.Range("E20").Value = "  +0.37%  "

$ws.Range("D21").Value = "423.27"
$ws.Range("E21").Value = "  -4.18%  "

$ws.Range("D22").Value = "97.07"
$ws.Range("E22").Value = "  +5.01%  "

$ws.Range("D23").Value = "3.47"
$ws.Range("E23").Value = "  -2.31%  "

$ws.Range("D24").Value = "4.12"
$ws.Range("E24").Value = "  +0.69%  "

$ws.Range("D25").Value = "14.17"
$ws.Range("E25").Value = "  -1.05%  "

$ws.Range("D26").Value = "10.95"
$ws.Range("E26").Value = "  -11.75%  "

$ws.Range("D27").Value = "10.59"
$ws.Range("E27").Value = "  -2.73%  "

$ws.Range("E28").Value = "  +1.85%  "

$ws.Range("D29").Value = "36.34"
$ws.Range("E29").Value = "  -1.50%  "

$ws.Range("D30").Value = "3.58"
$ws.Range("E30").Value = "  +23.92%  "

$ws.Range("D31").Value = "13.19"
$ws.Range("E31").Value = "  -3.17%  "

$ws.Range("B32").Value = "Bittensor"
$ws.Range("C32").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D32").Value = "673.85"
$ws.Range("E32").Value = "  -1.98%  "

$ws.Range("B33").Value = "Hedera"
$ws.Range("C33").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D33").Value = "0.127"
$ws.Range("E33").Value = "  -1.48%  "

$ws.Range("D34").Value = "6.80"
$ws.Range("E34").Value = "  -1.25%  "

$ws.Range("D35").Value = "65.40"
$ws.Range("E35").Value = "  -4.54%  "

$ws.Range("D36").Value = "42.03"
$ws.Range("E36").Value = "  +2.74%  "

$ws.Range("D37").Value = "0.419"
$ws.Range("E37").Value = "  -5.66%  "

$ws.Range("B38").Value = "PEPE"
$ws.Range("C38").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D38").Value = "0.0₃0828"
$ws.Range("E38").Value = "  -7.94%  "

$ws.Range("B39").Value = "Kaspa"
$ws.Range("C39").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D39").Value = "0.150"
$ws.Range("E39").Value = "  -0.81%  "

$ws.Range("D40").Value = "3.45"
$ws.Range("E40").Value = "  -1.91%  "

$ws.Range("D41").Value = "0.999"
$ws.Range("E41").Value = "  -0.01%  "

$ws.Range("B42").Value = "FirstDigitalUSD"
$ws.Range("C42").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D42").Value = "1.00"
$ws.Range("E42").Value = "  +0.14%  "

$ws.Range("B43").Value = "WEMIXToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D43").Value = "3.28"
$ws.Range("E43").Value = "  +5.75%  "

$ws.Range("D44").Value = "0.0480"
$ws.Range("E44").Value = "  -2.59%  "

$ws.Range("D45").Value = "0.149"
$ws.Range("E45").Value = "  +2.92%  "

$ws.Range("D46").Value = "9.49"
$ws.Range("E46").Value = "  +2.23%  "

$ws.Range("D47").Value = "2.57"
$ws.Range("E47").Value = "  -10.81%  "

$ws.Range("D48").Value = "3.30"
$ws.Range("E48").Value = "  -6.37%  "

$ws.Range("D49").Value = "2.98"
$ws.Range("E49").Value = "  -8.56%  "

$ws.Range("B50").Value = "FLOKI"
$ws.Range("C50").Value = "https://coinranking.com/coin/fmHk13Rqw+floki-floki"
$ws.Range("D50").Value = "0.000269"
$ws.Range("E50").Value = "  -5.24%  "

$ws.Range("B51").Value = "LidoDAOToken"
$ws.Range("C51").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D51").Value = "3.25"
$ws.Range("E51").Value = "  -4.66%  "
